$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "DBA" row (row 8)
$ws.Range("A8").Value = "DBA"
$ws.Range("B8").Value = 70
$ws.Range("C8").Value = 1
$ws.Range("D8").Formula = "=B8*168*C8"

# Add "Analista de Infra" row (row 9)
$ws.Range("A9").Value = "Analista de Infra"
$ws.Range("B9").Value = 37
$ws.Range("C9").Value = 1
$ws.Range("D9").Formula = "=B9*168*C9"

# Match the currency style used by the rest of column D (D2:D7)
$ws.Range("D8:D9").NumberFormat = $ws.Range("D7").NumberFormat
$ws.Range("D8:D9").Font.Color = $ws.Range("D7").Font.Color

# Move the formatted "total" cell down from D11 to D13
$ws.Range("D11").Copy($ws.Range("D13"))
$ws.Range("D11").ClearContents()

$null = $ws.Range("G9").Select()
